$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data row describing the MCH231-1 series (identifier, title,
# levelOfDescription, extentAndMedium, notes)
$ws.Range("A2").Value = "MCH231-1"
$ws.Range("C2").Value = "MWT PAPERS"
$ws.Range("E2").Value = "Series"
$ws.Range("F2").Value = "1 Box"
$ws.Range("G2").Value = "LOCATION: 24G | GRAP COUNT NUMER: NONE"

# Keep the default row heights explicit, as Excel records on save
$ws.Rows.Item(1).RowHeight = 15.75
$ws.Rows.Item(2).RowHeight = 15.75

# Re-apply the frozen header row, then select the newly entered row
$ws.Range("A2").Select() | Out-Null
[void]($excel.ActiveWindow.FreezePanes = $true)
$ws.Range("A2:I2").Select() | Out-Null
